$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-CellText "D2" '29.390.22'
Set-CellText "E2" '  +0.12%  '
Set-CellText "D3" '1.880.86'
Set-CellText "E3" '  +0.08%  '
Set-CellText "E4" '  +0.06%  '
Set-CellText "D5" '0.7167'
Set-CellText "E5" '  +0.72%  '
Set-CellText "D6" '243.54'
Set-CellText "E6" '  +0.54%  '
Set-CellText "D7" '1.002'
Set-CellText "E7" '  +0.10%  '
Set-CellText "D8" '0.07949'
Set-CellText "E8" '  -1.11%  '
Set-CellText "D9" '0.3139'
Set-CellText "E9" '  +0.28%  '
Set-CellText "D10" '24.85'
Set-CellText "E10" '  -1.58%  '
Set-CellText "E11" '  -3.47%  '
Set-CellText "D12" '1.945.58'
Set-CellText "E12" '  +4.02%  '
Set-CellText "D13" '94.62'
Set-CellText "E13" '  +3.54%  '
Set-CellText "D14" '5.202'
Set-CellText "E14" '  -1.06%  '
Set-CellText "D15" '0.7073'
Set-CellText "E15" '  -1.72%  '
Set-CellText "D16" '6.374'
Set-CellText "E16" '  +1.88%  '
Set-CellText "D17" '0.000008404'
Set-CellText "E17" '  +0.42%  '
Set-CellText "D18" '29.471.60'
Set-CellText "E18" '  +0.40%  '
Set-CellText "D19" '252.36'
Set-CellText "E19" '  +4.76%  '
Set-CellText "D20" '13.33'
Set-CellText "E20" '  +0.53%  '
Set-CellText "D21" '2.135.22'
Set-CellText "E21" '  +0.36%  '
Set-CellText "D23" '7.691'
Set-CellText "E23" '  -1.46%  '
Set-CellText "E24" '  +0.14%  '
Set-CellText "D25" '0.1576'
Set-CellText "E25" '  -0.95%  '
Set-CellText "D26" '9.052'
Set-CellText "E26" '  -0.06%  '
Set-CellText "D27" '161.57'
Set-CellText "E27" '  -0.89%  '
Set-CellText "D28" '18.92'
Set-CellText "E28" '  +1.89%  '
Set-CellText "D29" '1.509'
Set-CellText "E29" '  +0.16%  '
Set-CellText "D30" '4.408'
Set-CellText "E30" '  -0.31%  '
Set-CellText "D31" '4.312'
Set-CellText "E31" '  -0.62%  '
Set-CellText "D32" '1.235'
Set-CellText "E32" '  +3.90%  '
Set-CellText "D33" '0.05299'
Set-CellText "E33" '  -1.44%  '
Set-CellText "D34" '1.935'
Set-CellText "E34" '  -0.74%  '
Set-CellText "E35" '  +0.57%  '
Set-CellText "D36" '1.174'
Set-CellText "E36" '  -0.44%  '
Set-CellText "D38" '1.290.15'
Set-CellText "E38" '  -0.49%  '
Set-CellText "D39" '0.01880'
Set-CellText "D40" '2.767'
Set-CellText "E40" '  +1.03%  '
Set-CellText "D41" '6.398'
Set-CellText "E41" '  -2.79%  '
Set-CellText "D42" '0.9070'
Set-CellText "E42" '  +1.80%  '
Set-CellText "B43" 'Quant'
Set-CellText "C43" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-CellText "D43" '111.67'
Set-CellText "E43" '  +0.86%  '
Set-CellText "B44" 'Aave'
Set-CellText "C44" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-CellText "D44" '74.12'
Set-CellText "E44" '  +1.15%  '
Set-CellText "B45" 'PaxDollar'
Set-CellText "C45" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-CellText "D45" '1.001'
Set-CellText "E45" '  +0.07%  '
Set-CellText "B46" 'BabyDogeCoin'
Set-CellText "C46" 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-CellText "D46" '0.00000000130'
Set-CellText "E46" '  -0.45%  '
Set-CellText "D47" '2.039.13'
Set-CellText "E47" '  +0.67%  '
Set-CellText "E48" '  +0.41%  '
Set-CellText "D49" '0.5209'
Set-CellText "E49" '  -0.02%  '
Set-CellText "D50" '9.506'
Set-CellText "E50" '  +0.37%  '
Set-CellText "D51" '0.4341'
Set-CellText "E51" '  -0.53%  '
